$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

Set-TextValue "E2" "0.47%"
Set-TextValue "D3" "29.81"
Set-TextValue "E3" "9.54%"
Set-TextValue "D4" "5.173"
Set-TextValue "E4" "1.11%"
Set-TextValue "E5" "0.21%"
Set-TextValue "D6" "6.594"
Set-TextValue "E6" "1.57%"
Set-TextValue "B7" "GateToken"
Set-TextValue "C7" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D7" "3.098"
Set-TextValue "E7" "3.00%"
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.8560"
Set-TextValue "E8" "4.54%"
Set-TextValue "B9" "FTXToken"
Set-TextValue "C9" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D9" "0.8679"
Set-TextValue "E9" "1.88%"
Set-TextValue "B10" "One"
Set-TextValue "C10" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.01025"
Set-TextValue "E10" "1,606.59%"
Set-TextValue "B11" "WazirX"
Set-TextValue "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1368"
Set-TextValue "E11" "3.02%"
Set-TextValue "B12" "MandalaExchangeToken"
Set-TextValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07073"
Set-TextValue "E12" "2.00%"
Set-TextValue "B13" "BitrueCoin"
Set-TextValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.02930"
Set-TextValue "E13" "3.66%"
Set-TextValue "B14" "BitMartToken"
Set-TextValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09384"
Set-TextValue "E14" "-0.12%"
Set-TextValue "B15" "BitForexToken"
Set-TextValue "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001510"
Set-TextValue "E15" "-0.80%"
Set-TextValue "B16" "CoinExToken"
Set-TextValue "C16" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D16" "0.04179"
Set-TextValue "E16" "3.39%"
Set-TextValue "B17" "TigerCash"
Set-TextValue "C17" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D17" "0.006160"
Set-TextValue "E17" "-0.67%"
Set-TextValue "B18" "UpBots"
Set-TextValue "C18" "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D18" "0.007489"
Set-TextValue "E18" "5,073.35%"
Set-TextValue "B19" "LEO"
Set-TextValue "C19" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D19" "3.488"
Set-TextValue "E19" "-0.68%"
Set-TextValue "D20" "2.278"
Set-TextValue "E20" "-1.34%"
Set-TextValue "D22" "0.03444"
Set-TextValue "E22" "6.94%"
Set-TextValue "D23" "0.1311"
Set-TextValue "E23" "0.60%"
Set-TextValue "D24" "3.464"
Set-TextValue "E24" "-2.86%"
Set-TextValue "E25" "0.51%"
Set-TextValue "D26" "0.005015"
Set-TextValue "E26" "12.08%"
Set-TextValue "D27" "0.001220"
Set-TextValue "E27" "0.30%"
Set-TextValue "E28" "22.31%"
Set-TextValue "D40" "0.03755"
Set-TextValue "E40" "0.78%"
Set-TextValue "D41" "0.005759"
Set-TextValue "E41" "-1.64%"
Set-TextValue "D42" "0.1070"
Set-TextValue "E42" "1.07%"
Set-TextValue "D43" "0.002100"
Set-TextValue "E43" "-8.64%"
Set-TextValue "D44" "0.009576"
Set-TextValue "E44" "1.84%"
Set-TextValue "D45" "0.00005220"
Set-TextValue "E45" "1.36%"
Set-TextValue "E46" "0.07%"
Set-TextValue "E47" "-36.21%"
Set-TextValue "D48" "0.002531"
Set-TextValue "E48" "0.97%"
Set-TextValue "E49" "0.07%"
Set-TextValue "E50" "0.07%"

Write-Output "Applied 92 cell updates"
